# Refresh market-price derived columns (H:N) across all Leve tables.
# Source: scheduled market-data runner update (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1906.25
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1906.25
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1906.25
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = -2558.25
# Row 33
$ws.Range("H33").Value = 186.57895
$ws.Range("I33").Value = 222.5
$ws.Range("K33").Value = 222.5
$ws.Range("M33").Value = 6.5
# Row 80
$ws.Range("H80").Value = 1169.625
$ws.Range("I80").Value = 1456.4546
$ws.Range("J80").Value = 538.6
$ws.Range("K80").Value = 4369.3638
$ws.Range("L80").Value = 1615.8
$ws.Range("M80").Value = -3371.3638
$ws.Range("N80").Value = -3611.8
# Row 81
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
# Row 83
$ws.Range("H83").Value = 1169.625
$ws.Range("I83").Value = 1456.4546
$ws.Range("J83").Value = 538.6
$ws.Range("K83").Value = 13108.0914
$ws.Range("L83").Value = 4847.400000000001
$ws.Range("M83").Value = -8116.091400000001
$ws.Range("N83").Value = -14831.4
# Row 84
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
# Row 112
$ws.Range("H112").Value = 3430.28
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3430.28
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 10290.84
$ws.Range("M112").Value = $null
$ws.Range("N112").Value = -12506.84
# Row 141
$ws.Range("H141").Value = 3503410
$ws.Range("I141").Value = 5602056
$ws.Range("J141").Value = 5666.6665
$ws.Range("K141").Value = 16806168
$ws.Range("L141").Value = 16999.9995
$ws.Range("M141").Value = -16800988
$ws.Range("N141").Value = -27359.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 3810.3
$ws.Range("I74").Value = 3723.8235
$ws.Range("J74").Value = 4300.3335
$ws.Range("K74").Value = 3723.8235
$ws.Range("L74").Value = 4300.3335
$ws.Range("M74").Value = -2849.8235
$ws.Range("N74").Value = -6048.3335
# Row 77
$ws.Range("H77").Value = 3810.3
$ws.Range("I77").Value = 3723.8235
$ws.Range("J77").Value = 4300.3335
$ws.Range("K77").Value = 18619.1175
$ws.Range("L77").Value = 21501.6675
$ws.Range("M77").Value = -14251.1175
$ws.Range("N77").Value = -30237.6675
# Row 96
$ws.Range("H96").Value = 49972
$ws.Range("J96").Value = 49972
$ws.Range("L96").Value = 49972
$ws.Range("N96").Value = -55464
# Row 102
$ws.Range("H102").Value = 1045.4
$ws.Range("I102").Value = 998.25
$ws.Range("K102").Value = 998.25
$ws.Range("M102").Value = 623.75
# Row 122
$ws.Range("H122").Value = 1675
$ws.Range("I122").Value = 1661.3334
$ws.Range("K122").Value = 4984.0002
$ws.Range("M122").Value = -2534.0002
# Row 132
$ws.Range("H132").Value = 2203.5173
$ws.Range("I132").Value = 1370.2667
$ws.Range("K132").Value = 4110.800099999999
$ws.Range("M132").Value = -1580.800099999999

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 799
$ws.Range("I5").Value = 799
$ws.Range("K5").Value = 799
$ws.Range("M5").Value = -686
# Row 18
$ws.Range("H18").Value = 80011
$ws.Range("J18").Value = 80011
$ws.Range("L18").Value = 80011
$ws.Range("N18").Value = -81069
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
# Row 134
$ws.Range("H134").Value = 7578.7646
$ws.Range("I134").Value = 8432.24
$ws.Range("K134").Value = 25296.72
$ws.Range("M134").Value = -22761.72
# Row 141
$ws.Range("H141").Value = 41833.332
$ws.Range("J141").Value = 41833.332
$ws.Range("L141").Value = 41833.332
$ws.Range("N141").Value = -52193.332

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1350.7778
$ws.Range("J58").Value = 1699.5
$ws.Range("L58").Value = 1699.5
$ws.Range("N58").Value = -2105.5
# Row 86
$ws.Range("H86").Value = 1838.5
$ws.Range("I86").Value = 1838.5
$ws.Range("K86").Value = 1838.5
$ws.Range("M86").Value = -715.5
# Row 89
$ws.Range("H89").Value = 1838.5
$ws.Range("I89").Value = 1838.5
$ws.Range("K89").Value = 9192.5
$ws.Range("M89").Value = -3576.5
# Row 93
$ws.Range("H93").Value = 4933
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null
# Row 106
$ws.Range("H106").Value = 44850
$ws.Range("J106").Value = 44850
$ws.Range("L106").Value = 44850
$ws.Range("N106").Value = -47374
# Row 132
$ws.Range("H132").Value = 2234.6562
$ws.Range("I132").Value = 1183.3684
$ws.Range("K132").Value = 3550.1052
$ws.Range("M132").Value = -1020.1052
# Row 134
$ws.Range("H134").Value = 983.44446
$ws.Range("I134").Value = 944
$ws.Range("K134").Value = 2832
$ws.Range("M134").Value = -297
# Row 136
$ws.Range("H136").Value = 1350.7778
$ws.Range("J136").Value = 1699.5
$ws.Range("L136").Value = 5098.5
$ws.Range("N136").Value = -10198.5

$ws = $wb.Worksheets.Item("CUL")
# Row 130
$ws.Range("H130").Value = 1371.6
$ws.Range("I130").Value = 1089.5
$ws.Range("K130").Value = 3268.5
$ws.Range("M130").Value = 1751.5
# Row 131
$ws.Range("H131").Value = 21770150
$ws.Range("I131").Value = 71429110
$ws.Range("J131").Value = 44351.5
$ws.Range("K131").Value = 214287330
$ws.Range("L131").Value = 133054.5
$ws.Range("M131").Value = -214282290
$ws.Range("N131").Value = -143134.5

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 689.1875
$ws.Range("J97").Value = 1721
$ws.Range("L97").Value = 1721
$ws.Range("N97").Value = -2713
# Row 141
$ws.Range("H141").Value = 23000
$ws.Range("J141").Value = 23000
$ws.Range("L141").Value = 23000
$ws.Range("N141").Value = -33360

$ws = $wb.Worksheets.Item("LTW")
# Row 94
$ws.Range("H94").Value = 43540
$ws.Range("J94").Value = 43540
$ws.Range("L94").Value = 43540
$ws.Range("N94").Value = -44892
# Row 132
$ws.Range("H132").Value = 1996.3478
$ws.Range("I132").Value = 1599.8334
$ws.Range("J132").Value = 2136.2942
$ws.Range("K132").Value = 4799.5002
$ws.Range("L132").Value = 6408.882599999999
$ws.Range("M132").Value = -2269.5002
$ws.Range("N132").Value = -11468.8826

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 14285.714
$ws.Range("J18").Value = 14285.714
$ws.Range("L18").Value = 14285.714
$ws.Range("N18").Value = -14631.714
# Row 132
$ws.Range("H132").Value = 3710.111
$ws.Range("I132").Value = 1634.3334
$ws.Range("J132").Value = 4748
$ws.Range("K132").Value = 4903.0002
$ws.Range("L132").Value = 14244
$ws.Range("M132").Value = -2373.0002
$ws.Range("N132").Value = -19304
# Row 136
$ws.Range("H136").Value = 4093.2856
$ws.Range("I136").Value = 3941.0908
$ws.Range("K136").Value = 11823.2724
$ws.Range("M136").Value = -9273.2724

